{"js": "// The document body is a single paragraph holding the date, followed\n// by one table (20 rows x 5 cols) of arithmetic expressions. Office.js's\n// `body.paragraphs` collection walks every paragraph in the body in\n// document order, recursing into table cells, so index 0 is the date\n// paragraph and indices 1..100 are the 100 table-cell paragraphs\n// (row-major, left to right) -- the same order the diff lists them in.\nconst newValues = [\"2025-09-14 Sunday\", \"30+45=\", \"42-39=\", \"34-6=\", \"94-57=\", \"68-24=\", \"44+0=\", \"11+0=\", \"93-39=\", \"63-30=\", \"23+67=\", \"80-50=\", \"61-44=\", \"62-18=\", \"73-5=\", \"31+31=\", \"33+11=\", \"73-63=\", \"35+41=\", \"40-9=\", \"31-4=\", \"53-10=\", \"93-27=\", \"91-28=\", \"96-93=\", \"4+8=\", \"1+65=\", \"37+25=\", \"15-9=\", \"31+41=\", \"55-27=\", \"61-35=\", \"62+13=\", \"38+30=\", \"71+26=\", \"12+58=\", \"66+24=\", \"48+40=\", \"93-44=\", \"69-49=\", \"60-56=\", \"31-4=\", \"30+12=\", \"36+35=\", \"5+24=\", \"16+1=\", \"67-28=\", \"64+20=\", \"53-9=\", \"55+5=\", \"90-84=\", \"7+51=\", \"27-9=\", \"63+35=\", \"6+33=\", \"24-3=\", \"4+73=\", \"91-48=\", \"74-70=\", \"88-84=\", \"16+48=\", \"45-35=\", \"28+60=\", \"60+12=\", \"50-16=\", \"53+4=\", \"12+22=\", \"91-4=\", \"93-89=\", \"6+61=\", \"82-33=\", \"91-22=\", \"51+21=\", \"75-32=\", \"66+19=\", \"85-9=\", \"61-14=\", \"59-31=\", \"7+14=\", \"68-32=\", \"21+5=\", \"85+6=\", \"59-58=\", \"75-67=\", \"49+22=\", \"62+5=\", \"77-0=\", \"32-4=\", \"64-35=\", \"45-14=\", \"18+0=\", \"66+11=\", \"68-66=\", \"63-15=\", \"47-38=\", \"46-44=\", \"67-20=\", \"60-22=\", \"24-18=\", \"10+10=\", \"49-0=\"];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length && i < newValues.length; i++) {\n  items[i].getRange(\"Whole\").insertText(newValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document body is a single paragraph holding the date, followed by\n# one table (20 rows x 5 cols) of arithmetic expressions. Word's COM\n# Paragraphs collection walks every paragraph in the body in document\n# order, including one paragraph per table cell, plus an extra \"empty\"\n# paragraph at the end of every table row (the row/cell-mark, Chr 7)\n# and a final trailing paragraph after the table. Skip those empty\n# marker paragraphs (their Range.Text is nothing but the paragraph\n# mark / cell mark) and replace the text of the remaining 101 \"real\"\n# paragraphs -- the date, then the 100 cells row-major/left-to-right,\n# the same order the diff lists them in -- with the corresponding new\n# value. Assigning Range.Text replaces only the visible content and\n# keeps the paragraph's existing run formatting (font/size) intact.\n$newValues = @(\n    \"2025-09-14 Sunday\",\n    \"30+45=\",\n    \"42-39=\",\n    \"34-6=\",\n    \"94-57=\",\n    \"68-24=\",\n    \"44+0=\",\n    \"11+0=\",\n    \"93-39=\",\n    \"63-30=\",\n    \"23+67=\",\n    \"80-50=\",\n    \"61-44=\",\n    \"62-18=\",\n    \"73-5=\",\n    \"31+31=\",\n    \"33+11=\",\n    \"73-63=\",\n    \"35+41=\",\n    \"40-9=\",\n    \"31-4=\",\n    \"53-10=\",\n    \"93-27=\",\n    \"91-28=\",\n    \"96-93=\",\n    \"4+8=\",\n    \"1+65=\",\n    \"37+25=\",\n    \"15-9=\",\n    \"31+41=\",\n    \"55-27=\",\n    \"61-35=\",\n    \"62+13=\",\n    \"38+30=\",\n    \"71+26=\",\n    \"12+58=\",\n    \"66+24=\",\n    \"48+40=\",\n    \"93-44=\",\n    \"69-49=\",\n    \"60-56=\",\n    \"31-4=\",\n    \"30+12=\",\n    \"36+35=\",\n    \"5+24=\",\n    \"16+1=\",\n    \"67-28=\",\n    \"64+20=\",\n    \"53-9=\",\n    \"55+5=\",\n    \"90-84=\",\n    \"7+51=\",\n    \"27-9=\",\n    \"63+35=\",\n    \"6+33=\",\n    \"24-3=\",\n    \"4+73=\",\n    \"91-48=\",\n    \"74-70=\",\n    \"88-84=\",\n    \"16+48=\",\n    \"45-35=\",\n    \"28+60=\",\n    \"60+12=\",\n    \"50-16=\",\n    \"53+4=\",\n    \"12+22=\",\n    \"91-4=\",\n    \"93-89=\",\n    \"6+61=\",\n    \"82-33=\",\n    \"91-22=\",\n    \"51+21=\",\n    \"75-32=\",\n    \"66+19=\",\n    \"85-9=\",\n    \"61-14=\",\n    \"59-31=\",\n    \"7+14=\",\n    \"68-32=\",\n    \"21+5=\",\n    \"85+6=\",\n    \"59-58=\",\n    \"75-67=\",\n    \"49+22=\",\n    \"62+5=\",\n    \"77-0=\",\n    \"32-4=\",\n    \"64-35=\",\n    \"45-14=\",\n    \"18+0=\",\n    \"66+11=\",\n    \"68-66=\",\n    \"63-15=\",\n    \"47-38=\",\n    \"46-44=\",\n    \"67-20=\",\n    \"60-22=\",\n    \"24-18=\",\n    \"10+10=\",\n    \"49-0=\"\n)\n\n$d = $word.ActiveDocument\n\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    if ($i -ge $newValues.Count) {\n        break\n    }\n    $r = $p.Range\n    $plain = $r.Text.TrimEnd([char]13, [char]7)\n    if ($plain -ne \"\") {\n        $r.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
